# Added 4wk low sales check
# Update "Forecast Comparison" sheet metrics (columns H, I, J, L) for rows 2-17
$wb = $excel.ActiveWorkbook
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")

# Row 2 (W10)
$wsForecast.Range("H2").Value = 2.5
$wsForecast.Range("L2").Value = 0.95

# Row 3 (W11)
$wsForecast.Range("H3").Value = 1.5
$wsForecast.Range("L3").Value = 1.06

# Row 4 (W12)
$wsForecast.Range("H4").Value = 0.33
$wsForecast.Range("I4").Value = "High"
$wsForecast.Range("J4").Value = "Urgent"
$wsForecast.Range("L4").Value = 0.83

# Row 5 (W13)
$wsForecast.Range("H5").Value = 0
$wsForecast.Range("I5").Value = "High"
$wsForecast.Range("L5").Value = 1.03

# Row 6 (W14)
$wsForecast.Range("L6").Value = 1.13

# Row 7 (W15)
$wsForecast.Range("L7").Value = 1.07

# Row 8 (W16)
$wsForecast.Range("L8").Value = 0.93

# Row 9 (W17)
$wsForecast.Range("L9").Value = 0.98

# Row 10 (W18)
$wsForecast.Range("L10").Value = 1.07

# Row 11 (W19)
$wsForecast.Range("L11").Value = 0.92

# Row 12 (W20)
$wsForecast.Range("L12").Value = 0.96

# Row 13 (W21)
$wsForecast.Range("L13").Value = 0.87

# Row 14 (W22)
$wsForecast.Range("L14").Value = 0.95

# Row 15 (W23)
$wsForecast.Range("L15").Value = 0.93

# Row 16 (W24)
$wsForecast.Range("L16").Value = 1.02

# Row 17 (W25)
$wsForecast.Range("L17").Value = 1.18

# Update "Summary" sheet totals for 16-week and 4-week forecasts
# (force text number format so values stay stored as text, matching source data)
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B9").NumberFormat = "@"
$wsSummary.Range("B9").Value = "5"
$wsSummary.Range("B11").NumberFormat = "@"
$wsSummary.Range("B11").Value = "2"
